# Apply scheduled runner updates to Leve profit calculations across sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H92").Value = 995.6667
$ws.Range("I92").Value = 493.5
$ws.Range("K92").Value = 493.5
$ws.Range("M92").Value = 754.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H97").Value = 566.5
$ws.Range("J97").Value = 566.5
$ws.Range("L97").Value = 1699.5
$ws.Range("N97").Value = -2691.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 455.92856
$ws.Range("I98").Value = 390.33334
$ws.Range("J98").Value = 849.5
$ws.Range("K98").Value = 390.33334
$ws.Range("L98").Value = 849.5
$ws.Range("M98").Value = 1107.66666
$ws.Range("N98").Value = -3845.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H111").Value = 1133.1428
$ws.Range("I111").Value = 616.4
$ws.Range("J111").Value = 2425.0
$ws.Range("K111").Value = 1849.2
$ws.Range("L111").Value = 7275.0
$ws.Range("M111").Value = 1217.8
$ws.Range("N111").Value = -13409.0

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 10005.0
$ws.Range("I113").Value = 10005.0
$ws.Range("J113").Value = 0.0
$ws.Range("K113").Value = 10005.0
$ws.Range("L113").Value = 0.0
$ws.Range("M113").Value = -6751.0
$ws.Range("N113").ClearContents()

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H122").Value = 455.92856
$ws.Range("I122").Value = 390.33334
$ws.Range("J122").Value = 849.5
$ws.Range("K122").Value = 1171.00002
$ws.Range("L122").Value = 2548.5
$ws.Range("M122").Value = 1278.99998
$ws.Range("N122").Value = -7448.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H131").Value = 1565.0
$ws.Range("I131").Value = 1565.0
$ws.Range("K131").Value = 4695.0
$ws.Range("M131").Value = 345.0

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1158.2609
$ws.Range("I32").Value = 1078.0952
$ws.Range("K32").Value = 1078.0952
$ws.Range("M32").Value = -791.0952

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 4772.375
$ws.Range("I102").Value = 2163.3333
$ws.Range("J102").Value = 6337.8
$ws.Range("K102").Value = 2163.3333
$ws.Range("L102").Value = 6337.8
$ws.Range("M102").Value = -541.3332999999998
$ws.Range("N102").Value = -9581.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 0.0
$ws.Range("I132").Value = 0.0
$ws.Range("K132").Value = 0.0
$ws.Range("M132").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 7166.5
$ws.Range("J86").Value = 8199.8
$ws.Range("L86").Value = 8199.8
$ws.Range("N86").Value = -10445.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 7166.5
$ws.Range("J89").Value = 8199.8
$ws.Range("L89").Value = 40999.0
$ws.Range("N89").Value = -52231.0

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 862.0
$ws.Range("I16").Value = 484.4
$ws.Range("K16").Value = 484.4
$ws.Range("M16").Value = -197.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 815.9167
$ws.Range("I22").Value = 849.6
$ws.Range("J22").Value = 647.5
$ws.Range("K22").Value = 849.6
$ws.Range("L22").Value = 647.5
$ws.Range("M22").Value = -499.6
$ws.Range("N22").Value = -1347.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H68").Value = 50660.0
$ws.Range("J68").Value = 50660.0
$ws.Range("L68").Value = 50660.0
$ws.Range("N68").Value = -52158.0

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H71").Value = 50660.0
$ws.Range("J71").Value = 50660.0
$ws.Range("L71").Value = 151980.0
$ws.Range("N71").Value = -159468.0

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 3376.25
$ws.Range("J86").Value = 3502.5
$ws.Range("L86").Value = 3502.5
$ws.Range("N86").Value = -5748.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H89").Value = 3376.25
$ws.Range("J89").Value = 3502.5
$ws.Range("L89").Value = 17512.5
$ws.Range("N89").Value = -28744.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H113").Value = 862.0
$ws.Range("I113").Value = 484.4
$ws.Range("K113").Value = 484.4
$ws.Range("M113").Value = 1685.6

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 1896.2858
$ws.Range("I132").Value = 1896.2858
$ws.Range("K132").Value = 5688.857400000001
$ws.Range("M132").Value = -3158.857400000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H141").Value = 8250.0
$ws.Range("I141").Value = 1500.0
$ws.Range("J141").Value = 15000.0
$ws.Range("K141").Value = 4500.0
$ws.Range("L141").Value = 45000.0
$ws.Range("M141").Value = 680.0
$ws.Range("N141").Value = -55360.0

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2076.25
$ws.Range("I102").Value = 1621.125
$ws.Range("J102").Value = 3896.75
$ws.Range("K102").Value = 1621.125
$ws.Range("L102").Value = 3896.75
$ws.Range("M102").Value = 0.875
$ws.Range("N102").Value = -7140.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 89.666664
$ws.Range("I107").Value = 89.666664
$ws.Range("K107").Value = 89.666664
$ws.Range("M107").Value = 1830.333336

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2679.2
$ws.Range("I122").Value = 2679.2
$ws.Range("K122").Value = 8037.599999999999
$ws.Range("M122").Value = -5587.599999999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 2686.3333
$ws.Range("I126").Value = 2030.5454
$ws.Range("K126").Value = 6091.6362
$ws.Range("M126").Value = -3621.6362

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 1000012.0
$ws.Range("I132").Value = 1000012.0
$ws.Range("K132").Value = 3000036.0
$ws.Range("M132").Value = -2997506.0

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 8142.25
$ws.Range("I7").Value = 7756.5
$ws.Range("K7").Value = 7756.5
$ws.Range("M7").Value = -7644.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 5306.125
$ws.Range("I82").Value = 700.0
$ws.Range("K82").Value = 700.0
$ws.Range("M82").Value = -339.0

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H85").Value = 5306.125
$ws.Range("I85").Value = 700.0
$ws.Range("K85").Value = 700.0
$ws.Range("M85").Value = 548.0

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 1098.7333
$ws.Range("I93").Value = 1116.6364
$ws.Range("K93").Value = 1116.6364
$ws.Range("M93").Value = 131.3635999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 2999.0
$ws.Range("I122").Value = 2999.0
$ws.Range("J122").Value = 0.0
$ws.Range("K122").Value = 8997.0
$ws.Range("L122").Value = 0.0
$ws.Range("M122").Value = -6547.0
$ws.Range("N122").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 8142.25
$ws.Range("I126").Value = 7756.5
$ws.Range("K126").Value = 23269.5
$ws.Range("M126").Value = -20799.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 351.0
$ws.Range("J81").Value = 0.0
$ws.Range("L81").Value = 0.0
$ws.Range("N81").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H84").Value = 351.0
$ws.Range("J84").Value = 0.0
$ws.Range("L84").Value = 0.0
$ws.Range("N84").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 6242.3335
$ws.Range("I126").Value = 5485.0
$ws.Range("J126").Value = 6999.6665
$ws.Range("K126").Value = 16455.0
$ws.Range("L126").Value = 20998.9995
$ws.Range("M126").Value = -13985.0
$ws.Range("N126").Value = -25938.9995
